$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 778, shifting existing rows 778-819 down to 779-820.
$ws.Rows.Item(778).Insert()

# Populate the newly inserted row 778 with the new data point.
# Use a leading apostrophe for the date so it is stored as literal text
# (matching the inlineStr date text used throughout column A) rather than
# being auto-converted into a date serial number.
$ws.Range("A778").Value = "'2026/02/06"
$ws.Range("A778").Style = "Normal"

$ws.Range("B778").Value = "金"
$ws.Range("C778").Value = 3
$ws.Range("D778").Value = 201
